# SB all response collected for pairwise
# Adds three new "SB" tracker rows (row 16 and 17 are brand new, row 15
# gains a value in column E) to the single worksheet of the tracker
# workbook, mirroring the same Date/expert/Task_type pattern already
# used by the other rows for Oct-1-2023 / SB.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (new): paiewise_resub / Oct-1-2023 / SB / re_submitted... / master_all_responses... ---
$ws.Range("A16").Value2 = "paiewise_resub"

# Copy B15 ("Oct-1-2023") as a value so the destination keeps the same
# shared-string text type instead of Excel's automatic date parsing.
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("C16").Value2 = "SB"
$ws.Range("E16").Value2 = "master_all_responses_SB_resub_Oct-01-2023.csv"
$ws.Range("D16").Value2 = "re_submitted_tracker_SB_Oct-01-2023.csv"

# --- Row 15 (existing): add the missing response_collected file name ---
$ws.Range("E15").Value2 = "master_worker_response_tracke_SB_Oct-01-2023.csv"

# --- Row 17 (new): pairwise_resub_video_not_play / Oct-1-2023 / SB / all_submitted_tracker... ---
$ws.Range("A17").Value2 = "pairwise_resub_video_not_play"

$ws.Range("B15").Copy()
$ws.Range("B17").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("C17").Value2 = "SB"
$ws.Range("D17").Value2 = "all_submitted_tracker_SB_video_no_play_Oct-01-2023.csv"

# Page was switched to portrait orientation for printing.
$ws.PageSetup.Orientation = 1

# Selection ends up on E18 (first empty row after the new data) after entry.
$ws.Range("E18").Select()
